$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2835.923
$ws.Range("I17").Value = 3813
$ws.Range("J17").Value = 1998.4286
$ws.Range("K17").Value = 11439
$ws.Range("L17").Value = 5995.2858
$ws.Range("M17").Value = -11271
$ws.Range("N17").Value = -6331.2858
# Row 98
$ws.Range("H98").Value = 1108.4286
$ws.Range("I98").Value = 1108.4286
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1108.4286
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = 389.5714
$ws.Range("M98").ClearContents()
# Row 122
$ws.Range("H122").Value = 1108.4286
$ws.Range("I122").Value = 1108.4286
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3325.2858
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -875.2857999999997
$ws.Range("M122").ClearContents()
# Row 131
$ws.Range("H131").Value = 4275.933
$ws.Range("I131").Value = 1404.3334
$ws.Range("J131").Value = 8583.333000000001
$ws.Range("K131").Value = 4213.0002
$ws.Range("L131").Value = 25749.999
$ws.Range("M131").Value = 826.9997999999996
$ws.Range("N131").Value = -35829.999

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 693.3333
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 790
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 790
$ws.Range("M4").Value = -384
$ws.Range("N4").Value = -1022
# Row 74
$ws.Range("H74").Value = 2797.4119
$ws.Range("I74").Value = 923.61536
$ws.Range("K74").Value = 923.61536
$ws.Range("M74").Value = -49.61536000000001
# Row 77
$ws.Range("H77").Value = 2797.4119
$ws.Range("I77").Value = 923.61536
$ws.Range("K77").Value = 4618.0768
$ws.Range("M77").Value = -250.0767999999998
# Row 132
$ws.Range("H132").Value = 2182.5557
$ws.Range("I132").Value = 1952.4
$ws.Range("K132").Value = 5857.200000000001
$ws.Range("M132").Value = -3327.200000000001

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 95999
$ws.Range("J52").Value = 95999
$ws.Range("L52").Value = 95999
$ws.Range("N52").Value = -96587
# Row 58
$ws.Range("H58").Value = 3050.7273
$ws.Range("I58").Value = 2008.1428
$ws.Range("J58").Value = 3537.2666
$ws.Range("K58").Value = 2008.1428
$ws.Range("L58").Value = 3537.2666
$ws.Range("M58").Value = -1805.1428
$ws.Range("N58").Value = -3943.2666
# Row 99
$ws.Range("H99").Value = 10764.936
$ws.Range("I99").Value = 6271.75
$ws.Range("K99").Value = 6271.75
$ws.Range("M99").Value = -4773.75
# Row 125
$ws.Range("H125").Value = 84500
$ws.Range("J125").Value = 84500
$ws.Range("L125").Value = 84500
$ws.Range("N125").Value = -89420
# Row 126
$ws.Range("H126").Value = 10764.936
$ws.Range("I126").Value = 6271.75
$ws.Range("K126").Value = 18815.25
$ws.Range("M126").Value = -16345.25
# Row 132
$ws.Range("H132").Value = 3662.375
$ws.Range("I132").Value = 2360.75
$ws.Range("J132").Value = 4964
$ws.Range("K132").Value = 7082.25
$ws.Range("L132").Value = 14892
$ws.Range("M132").Value = -4552.25
$ws.Range("N132").Value = -19952
# Row 136
$ws.Range("H136").Value = 3050.7273
$ws.Range("I136").Value = 2008.1428
$ws.Range("J136").Value = 3537.2666
$ws.Range("K136").Value = 6024.428400000001
$ws.Range("L136").Value = 10611.7998
$ws.Range("M136").Value = -3474.428400000001
$ws.Range("N136").Value = -15711.7998

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 100078.8
$ws.Range("I2").Value = 200025.8
$ws.Range("K2").Value = 1200154.8
$ws.Range("M2").Value = -1200041.8
# Row 50
$ws.Range("H50").Value = 211
$ws.Range("I50").Value = 156.8
$ws.Range("K50").Value = 470.4
$ws.Range("M50").Value = 10.59999999999997
# Row 53
$ws.Range("H53").Value = 211
$ws.Range("I53").Value = 156.8
$ws.Range("K53").Value = 470.4
$ws.Range("M53").Value = 10.59999999999997
# Row 80
$ws.Range("H80").Value = 3895.1667
$ws.Range("I80").Value = 2844.25
$ws.Range("K80").Value = 8532.75
$ws.Range("M80").Value = -7596.75
# Row 83
$ws.Range("H83").Value = 3895.1667
$ws.Range("I83").Value = 2844.25
$ws.Range("K83").Value = 25598.25
$ws.Range("M83").Value = -20918.25
# Row 86
$ws.Range("H86").Value = 274.5
$ws.Range("I86").Value = 199.5
$ws.Range("K86").Value = 598.5
$ws.Range("M86").Value = 587.5
# Row 89
$ws.Range("H89").Value = 274.5
$ws.Range("I89").Value = 199.5
$ws.Range("K89").Value = 1795.5
$ws.Range("M89").Value = 4132.5
# Row 92
$ws.Range("H92").Value = 658
$ws.Range("I92").Value = 649.5
$ws.Range("J92").Value = 675
$ws.Range("K92").Value = 1948.5
$ws.Range("L92").Value = 2025
$ws.Range("M92").Value = -700.5
$ws.Range("N92").Value = -4521
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
# Row 107
$ws.Range("H107").Value = 717.8484999999999
$ws.Range("J107").Value = 707.29034
$ws.Range("L107").Value = 2121.87102
$ws.Range("N107").Value = -5961.87102
# Row 113
$ws.Range("H113").Value = 1606.5834
$ws.Range("I113").Value = 2333
$ws.Range("J113").Value = 1364.4445
$ws.Range("K113").Value = 6999
$ws.Range("L113").Value = 4093.3335
$ws.Range("M113").Value = -4829
$ws.Range("N113").Value = -8433.333500000001

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 728.3
$ws.Range("I2").Value = 182.6
$ws.Range("J2").Value = 1274
$ws.Range("K2").Value = 182.6
$ws.Range("L2").Value = 1274
$ws.Range("M2").Value = -69.59999999999999
$ws.Range("N2").Value = -1500
# Row 102
$ws.Range("H102").Value = 1539.2759
$ws.Range("I102").Value = 364.09525
$ws.Range("J102").Value = 4624.125
$ws.Range("K102").Value = 364.09525
$ws.Range("L102").Value = 4624.125
$ws.Range("M102").Value = 1257.90475
$ws.Range("N102").Value = -7868.125
# Row 107
$ws.Range("H107").Value = 668.8
$ws.Range("I107").Value = 468.16666
$ws.Range("J107").Value = 854
$ws.Range("K107").Value = 468.16666
$ws.Range("L107").Value = 854
$ws.Range("M107").Value = 1451.83334
$ws.Range("N107").Value = -4694
# Row 122
$ws.Range("H122").Value = 70536.8
$ws.Range("I122").Value = 3867.3845
$ws.Range("K122").Value = 11602.1535
$ws.Range("M122").Value = -9152.1535
# Row 126
$ws.Range("H126").Value = 3875.5
$ws.Range("I126").Value = 3022.8572
$ws.Range("J126").Value = 4418.091
$ws.Range("K126").Value = 9068.571599999999
$ws.Range("L126").Value = 13254.273
$ws.Range("M126").Value = -6598.571599999999
$ws.Range("N126").Value = -18194.273

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2371.75
$ws.Range("J7").Value = 1996.5
$ws.Range("L7").Value = 1996.5
$ws.Range("N7").Value = -2220.5
# Row 40
$ws.Range("H40").Value = 1132.1
$ws.Range("I40").Value = 1039.5
$ws.Range("K40").Value = 1039.5
$ws.Range("M40").Value = -903.5
# Row 100
$ws.Range("H100").Value = 14000
$ws.Range("J100").Value = 14000
$ws.Range("L100").Value = 14000
$ws.Range("N100").Value = -15082
# Row 122
$ws.Range("H122").Value = 4608.85
$ws.Range("I122").Value = 2667.7
$ws.Range("K122").Value = 8003.099999999999
$ws.Range("M122").Value = -5553.099999999999
# Row 126
$ws.Range("H126").Value = 2371.75
$ws.Range("J126").Value = 1996.5
$ws.Range("L126").Value = 5989.5
$ws.Range("N126").Value = -10929.5
# Row 132
$ws.Range("H132").Value = 3726.9546
$ws.Range("I132").Value = 3360.2593
$ws.Range("K132").Value = 10080.7779
$ws.Range("M132").Value = -7550.777900000001

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1963.4736
$ws.Range("I122").Value = 887.4
$ws.Range("J122").Value = 5998.75
$ws.Range("K122").Value = 2662.2
$ws.Range("L122").Value = 17996.25
$ws.Range("M122").Value = -212.1999999999998
$ws.Range("N122").Value = -22896.25
